# Generate Report for Handoff
#
# Row 7 (previously "c756cc0b..." / In Translation) and Row 8 (previously
# "0a4156de..." / Ready for handoff) swap positions: the 0a4156de entry
# (already ready for handoff) now sorts ahead of c756cc0b, which itself
# has just been handed off (new, later timestamps).
#
# Applies to all three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet (columns: A=File Name, B=Path And Name, C=Extension,
# D=Publish URL, E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A7").Value = "0a4156de-673c-4b84-942e-2bd0c324cced.md"
$wsOverview.Range("B7").Value = "e2e\0a4156de-673c-4b84-942e-2bd0c324cced.md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-09-03 16:48:18"

$wsOverview.Range("A8").Value = "c756cc0b-d867-45c3-92e4-83d06984c87e.md"
$wsOverview.Range("B8").Value = "e2e\c756cc0b-d867-45c3-92e4-83d06984c87e.md"
$wsOverview.Range("E8").Value = "Ready for handoff"
$wsOverview.Range("F8").Value = "Ready for handoff"
$wsOverview.Range("G8").Value = "2016-09-03 16:51:17"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$7') {
        $h.TextToDisplay = "e2e\0a4156de-673c-4b84-942e-2bd0c324cced.md"
    } elseif ($addr -eq '$B$8') {
        $h.TextToDisplay = "e2e\c756cc0b-d867-45c3-92e4-83d06984c87e.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet (A=Source File Name, C=Status, G=Latest Handoff File,
# H=Latest Handoff Datetime)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A7").Value = "0a4156de-673c-4b84-942e-2bd0c324cced.md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("G7").Value = "0a4156de-673c-4b84-942e-2bd0c324cced.5fde648ae5e7fe0db87bfc415a0db69f46ad3b97.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-09-03 16:48:13"

$wsZhCn.Range("A8").Value = "c756cc0b-d867-45c3-92e4-83d06984c87e.md"
$wsZhCn.Range("C8").Value = "Ready for handoff"
$wsZhCn.Range("G8").Value = "c756cc0b-d867-45c3-92e4-83d06984c87e.33059936ab1f38d4316c666616f6699d71ea8d57.zh-cn.xlf"
$wsZhCn.Range("H8").Value = "2016-09-03 16:51:12"

foreach ($h in $wsZhCn.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$7') {
        $h.TextToDisplay = "0a4156de-673c-4b84-942e-2bd0c324cced.md"
    } elseif ($addr -eq '$A$8') {
        $h.TextToDisplay = "c756cc0b-d867-45c3-92e4-83d06984c87e.md"
    }
}

# ---------------------------------------------------------------------
# de-de sheet (A=Source File Name, C=Status, G=Latest Handoff File,
# H=Latest Handoff Datetime)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A7").Value = "0a4156de-673c-4b84-942e-2bd0c324cced.md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("G7").Value = "0a4156de-673c-4b84-942e-2bd0c324cced.5fde648ae5e7fe0db87bfc415a0db69f46ad3b97.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-09-03 16:48:18"

$wsDeDe.Range("A8").Value = "c756cc0b-d867-45c3-92e4-83d06984c87e.md"
$wsDeDe.Range("C8").Value = "Ready for handoff"
$wsDeDe.Range("G8").Value = "c756cc0b-d867-45c3-92e4-83d06984c87e.33059936ab1f38d4316c666616f6699d71ea8d57.de-de.xlf"
$wsDeDe.Range("H8").Value = "2016-09-03 16:51:17"

foreach ($h in $wsDeDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$7') {
        $h.TextToDisplay = "0a4156de-673c-4b84-942e-2bd0c324cced.md"
    } elseif ($addr -eq '$A$8') {
        $h.TextToDisplay = "c756cc0b-d867-45c3-92e4-83d06984c87e.md"
    }
}
